$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mis-encoded French accented canton names in column A
$ws.Range("A24").Value = "Bâle-Ville"
$ws.Range("A25").Value = "Bâle-Ville"
$ws.Range("A26").Value = "Bâle-Campagne"
$ws.Range("A27").Value = "Bâle-Campagne"
$ws.Range("A30").Value = "Appenzell Rhodes-Extérieures"
$ws.Range("A31").Value = "Appenzell Rhodes-Extérieures"
$ws.Range("A32").Value = "Appenzell Rhodes-Intérieures"
$ws.Range("A33").Value = "Appenzell Rhodes-Intérieures"
$ws.Range("A48").Value = "Neuchâtel"
$ws.Range("A49").Value = "Neuchâtel"
$ws.Range("A50").Value = "Genève"
$ws.Range("A51").Value = "Genève"

# Reset the view back to the top-left and set the active cell selection
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A49").Select()
